$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert two new columns before column D (shifts D:K -> F:M)
$ws.Range("D:E").Insert()

# Step 2: Copy number formats into the new D/E columns from F/G so they match styling
$ws.Range("F5:F102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Range("G5:G102").Copy()
$ws.Range("E5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: Populate the new D/E columns with the new quarter's data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 124100
$ws.Range("E8").Value = 127900
$ws.Range("D9").Value = 103000
$ws.Range("E9").Value = 102200
$ws.Range("D10").Value = 21100
$ws.Range("E10").Value = 25700
$ws.Range("D12").Value = 2100
$ws.Range("E12").Value = 1800
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 122500
$ws.Range("E17").Value = 117500
$ws.Range("D18").Value = 1600
$ws.Range("E18").Value = 10400
$ws.Range("D20").Value = -5800
$ws.Range("E20").Value = -5700
$ws.Range("D21").Value = -2400
$ws.Range("E21").Value = 6300
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = -4100
$ws.Range("E23").Value = 4700
$ws.Range("D24").Value = -1000
$ws.Range("E24").Value = 1200
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -3100
$ws.Range("E26").Value = 3500
$ws.Range("D27").Value = -3100
$ws.Range("E27").Value = 3500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 5800
$ws.Range("E32").Value = 5700
$ws.Range("D33").Value = -3100
$ws.Range("E33").Value = 3500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -3100
$ws.Range("E35").Value = 3500
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 8600
$ws.Range("E41").Value = 6200
$ws.Range("D42").Value = 2900
$ws.Range("E42").Value = 3200
$ws.Range("D43").Value = 94600
$ws.Range("E43").Value = 91900
$ws.Range("D44").Value = 215800
$ws.Range("E44").Value = 197400
$ws.Range("D45").Value = 7900
$ws.Range("E45").Value = 10200
$ws.Range("D46").Value = 329800
$ws.Range("E46").Value = 308900
$ws.Range("D47").Value = 223000
$ws.Range("E47").Value = 230400
$ws.Range("D48").Value = 32300
$ws.Range("E48").Value = 30500
$ws.Range("D49").Value = 12200
$ws.Range("E49").Value = 5900
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 8500
$ws.Range("E52").Value = 8200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 605800
$ws.Range("E54").Value = 584000
$ws.Range("D57").Value = 101700
$ws.Range("E57").Value = 92700
$ws.Range("D58").Value = 82100
$ws.Range("E58").Value = 56600
$ws.Range("D59").Value = 66900
$ws.Range("E59").Value = 67000
$ws.Range("D60").Value = 250700
$ws.Range("E60").Value = 216300
$ws.Range("D61").Value = 25100
$ws.Range("E61").Value = 26000
$ws.Range("D62").Value = 50300
$ws.Range("E62").Value = 59500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 326100
$ws.Range("E66").Value = 301800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 74200
$ws.Range("E72").Value = 77300
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 279800
$ws.Range("E76").Value = 282200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -3100
$ws.Range("E81").Value = 3500
$ws.Range("D83").Value = 1700
$ws.Range("E83").Value = 1600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -13900
$ws.Range("E89").Value = -5500
$ws.Range("D91").Value = -3300
$ws.Range("E91").Value = -3700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -7800
$ws.Range("E94").Value = -3800
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 24200
$ws.Range("E100").Value = 3200
$ws.Range("D101").Value = -100
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 2400
$ws.Range("E102").Value = -6100

# Step 4: Apply restated values to existing columns (financial data revisions)
$ws.Range("F8").Value = 91700
$ws.Range("H8").Value = 102900
$ws.Range("I8").Value = 110300
$ws.Range("J8").Value = 94700
$ws.Range("H9").Value = 76800
$ws.Range("I9").Value = 84200
$ws.Range("F10").Value = 16400
$ws.Range("H10").Value = 26100
$ws.Range("I10").Value = 26100
$ws.Range("J10").Value = 25900
$ws.Range("F17").Value = 93500
$ws.Range("H17").Value = 94500
$ws.Range("I17").Value = 97500
$ws.Range("J17").Value = 79100
$ws.Range("F18").Value = -1800
$ws.Range("H18").Value = 8400
$ws.Range("I18").Value = 12800
$ws.Range("J18").Value = 15600
$ws.Range("H20").Value = -4000
$ws.Range("I20").Value = -3600
$ws.Range("F21").Value = -5400
$ws.Range("H21").Value = 5600
$ws.Range("I21").Value = 10300
$ws.Range("J21").Value = 13300
$ws.Range("F23").Value = -6900
$ws.Range("H23").Value = 4500
$ws.Range("I23").Value = 9200
$ws.Range("J23").Value = 12300
$ws.Range("F24").Value = -1400
$ws.Range("H24").Value = 200
$ws.Range("I24").Value = 3600
$ws.Range("J24").Value = 4400
$ws.Range("F26").Value = -5500
$ws.Range("H26").Value = 4300
$ws.Range("I26").Value = 5600
$ws.Range("J26").Value = 7800
$ws.Range("F27").Value = -5500
$ws.Range("H27").Value = 4300
$ws.Range("I27").Value = 5600
$ws.Range("J27").Value = 7800
$ws.Range("F29").Value = "NA"
$ws.Range("I29").Value = "NA"
$ws.Range("J29").Value = "NA"
$ws.Range("K29").Value = "NA"
$ws.Range("L29").Value = "NA"
$ws.Range("M29").Value = "NA"
$ws.Range("H32").Value = 4000
$ws.Range("I32").Value = 3600
$ws.Range("F33").Value = -5500
$ws.Range("H33").Value = -2500
$ws.Range("I33").Value = 5600
$ws.Range("J33").Value = 7800
$ws.Range("F35").Value = -5500
$ws.Range("H35").Value = -2500
$ws.Range("I35").Value = 5600
$ws.Range("J35").Value = 7800
$ws.Range("G43").Value = 157500
$ws.Range("G44").Value = 329900
$ws.Range("G45").Value = 20100
$ws.Range("G46").Value = 287600
$ws.Range("G47").Value = 228300
$ws.Range("G49").Value = 10100
$ws.Range("G52").Value = 309200
$ws.Range("G54").Value = 552400
$ws.Range("G59").Value = 83500
$ws.Range("G60").Value = 197300
$ws.Range("G66").Value = 265500
$ws.Range("G72").Value = 78500
$ws.Range("G76").Value = 286900
$ws.Range("F81").Value = -5500
$ws.Range("H81").Value = -2500
$ws.Range("I81").Value = 5600
$ws.Range("J81").Value = 7800
$ws.Range("I91").Value = -1900
$ws.Range("J91").Value = -600
